$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Select cust_country group by  " + "cust_country"  -> single run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Select cust_country group by  cust_country",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Select cust_country group by  cust_country", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Select " + "min" + "(rollnumber) from student" -> single run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Select min(rollnumber) from student",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Select min(rollnumber) from student", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Select " + "max" + "(rollnumber) from student" -> single run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Select max(rollnumber) from student",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Select max(rollnumber) from student", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Select " + "avg" + "(rollnumber) from student" -> single run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Select avg(rollnumber) from student",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Select avg(rollnumber) from student", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "Select name from student order by  name ASC (or desc)" (double space)
#    -> "Select name from student order " / "by name" / " ASC (or desc)"
#    (single space between "order" and "name" in the result)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Select name from student order by  name ASC (or desc)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Select name from student order by name ASC (or desc)", 2) | Out-Null

# locate that paragraph again and split its single run into three runs.
# Trick: temporarily break the paragraph into three mini-paragraphs at the
# desired boundaries (a bare CR on a collapsed range splits cleanly), then
# delete the two inserted paragraph marks so the text rejoins into a single
# paragraph again -- but the runs that were on either side of each deleted
# mark stay distinct, with none of the empty-rPr residue a plain character
# format toggle would leave behind.
$pIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Select name from student order by name ASC (or desc)") {
        $pIndex = $i
        break
    }
}

$pStart = $d.Paragraphs($pIndex).Range.Start
$part1 = "Select name from student order "
$part2 = "by name"
$len1 = $part1.Length
$len2 = $part2.Length

# split at the later boundary first so the earlier offset stays valid
$posB = $pStart + $len1 + $len2
$d.Range($posB, $posB).Text = "`r"
$posA = $pStart + $len1
$d.Range($posA, $posA).Text = "`r"

# merge the three mini-paragraphs back into one, later mark first
$markEnd2 = $d.Paragraphs($pIndex + 1).Range.End
$d.Range($markEnd2 - 1, $markEnd2).Delete() | Out-Null
$markEnd1 = $d.Paragraphs($pIndex).Range.End
$d.Range($markEnd1 - 1, $markEnd1).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 6) Insert two new paragraphs ("Between, like, %%, _, " and a blank one)
#    right before the paragraph that holds the _GoBack bookmark, then give
#    that bookmark paragraph the GitHub link text, and finally delete the
#    blank paragraph that used to trail it.
# ---------------------------------------------------------------------------

# find the paragraph that contains the _GoBack bookmark by inspecting the
# raw paragraph XML (Range.Bookmarks.Count is unreliable in this host)
$bookmarkParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $xml = $d.Paragraphs($i).Range.WordOpenXML
    if ($xml -like "*_GoBack*") {
        $bookmarkParaIndex = $i
        break
    }
}

$bm = $d.Paragraphs($bookmarkParaIndex)

# insert a new blank paragraph right before the bookmark paragraph (a bare
# carriage return on a collapsed range splits cleanly, with no stray empty
# run left behind), fill it with the "Between, like, ..." text
$bmPos = $bm.Range.Start
$d.Range($bmPos, $bmPos).Text = "`r"
$betweenPara = $d.Paragraphs($bookmarkParaIndex)
$betweenPara.Range.Text = "Between, like, %%, _, "

# insert another new blank paragraph right before the (shifted) bookmark
# paragraph; leave it empty
$bm2 = $d.Paragraphs($bookmarkParaIndex + 1)
$bm2Pos = $bm2.Range.Start
$d.Range($bm2Pos, $bm2Pos).Text = "`r"

# the bookmark paragraph has shifted two slots down; give it the URL text,
# inserted right at its start (i.e. before the bookmark start/end marks)
$bmFinal = $d.Paragraphs($bookmarkParaIndex + 2)
$urlPoint = $d.Range($bmFinal.Range.Start, $bmFinal.Range.Start)
$urlPoint.InsertBefore("https://github.com/SanjayTamboli7/ST-SE-Tops/tree/main/Database/Assignments/Module%205")

# remove the paragraph that used to trail the bookmark paragraph: it is now
# the very last (empty) paragraph of the body, right after the bookmark one
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastText = $lastPara.Range.Text.TrimEnd([char]13, [char]7)
if ($lastText -eq "") {
    $delRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
    $delRange.Delete() | Out-Null
}

Write-Output "edits applied"
